$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1011-1012, shifting existing rows 1011..1096 down to 1013..1098.
$ws.Rows("1011:1012").Insert()

# Row 1011 - new data row
$ws.Cells.Item(1011, 1).Value = 3
$ws.Cells.Item(1011, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1011, 3).Value = "Coquimbo"
$ws.Cells.Item(1011, 4).Value = 44769
$ws.Cells.Item(1011, 5).Value = 5
$ws.Cells.Item(1011, 6).Value = "Fruta"
$ws.Cells.Item(1011, 7).Value = 100102
$ws.Cells.Item(1011, 8).Value = "Cítricos"
$ws.Cells.Item(1011, 9).Value = 100102003
$ws.Cells.Item(1011, 10).Value = "Limón"
$ws.Cells.Item(1011, 11).Value = "Sin especificar"
$ws.Cells.Item(1011, 12).Value = "1a amarillo"
$ws.Cells.Item(1011, 13).Value = 310
$ws.Cells.Item(1011, 14).Value = 2700
$ws.Cells.Item(1011, 15).Value = 3000
$ws.Cells.Item(1011, 16).Value = 2840
$ws.Cells.Item(1011, 17).Value = '$/malla 16 kilos'
$ws.Cells.Item(1011, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(1011, 19).Value = 178
$ws.Cells.Item(1011, 20).Value = 16

# Row 1012 - new data row
$ws.Cells.Item(1012, 1).Value = 3
$ws.Cells.Item(1012, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1012, 3).Value = "Coquimbo"
$ws.Cells.Item(1012, 4).Value = 44769
$ws.Cells.Item(1012, 5).Value = 5
$ws.Cells.Item(1012, 6).Value = "Fruta"
$ws.Cells.Item(1012, 7).Value = 100102
$ws.Cells.Item(1012, 8).Value = "Cítricos"
$ws.Cells.Item(1012, 9).Value = 100102003
$ws.Cells.Item(1012, 10).Value = "Limón"
$ws.Cells.Item(1012, 11).Value = "Sin especificar"
$ws.Cells.Item(1012, 12).Value = "2a amarillo"
$ws.Cells.Item(1012, 13).Value = 307
$ws.Cells.Item(1012, 14).Value = 2300
$ws.Cells.Item(1012, 15).Value = 2500
$ws.Cells.Item(1012, 16).Value = 2401
$ws.Cells.Item(1012, 17).Value = '$/malla 16 kilos'
$ws.Cells.Item(1012, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(1012, 19).Value = 150
$ws.Cells.Item(1012, 20).Value = 16
